$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D16").Value = 44592
$ws.Range("I16").Value = "Tercera"
$ws.Range("J16").Value = 200
$ws.Range("K16").Value = 1500
$ws.Range("L16").Value = 1800
$ws.Range("M16").Value = 1650
$ws.Range("P16").Value = 550
$ws.Range("D17").Value = 44537
$ws.Range("I17").Value = "Primera"
$ws.Range("J17").Value = 250
$ws.Range("K17").Value = 1400
$ws.Range("L17").Value = 1500
$ws.Range("M17").Value = 1450
$ws.Range("P17").Value = 483
$ws.Range("D18").Value = 44320
$ws.Range("I18").Value = "Primera"
$ws.Range("J18").Value = 200
$ws.Range("K18").Value = 1400
$ws.Range("L18").Value = 1500
$ws.Range("M18").Value = 1450
$ws.Range("P18").Value = 483
$ws.Range("D19").Value = 44320
$ws.Range("I19").Value = "Segunda"
$ws.Range("J19").Value = 200
$ws.Range("K19").Value = 1000
$ws.Range("L19").Value = 1200
$ws.Range("M19").Value = 1100
$ws.Range("P19").Value = 367
$ws.Range("D20").Value = 44383
$ws.Range("I20").Value = "Segunda"
$ws.Range("J20").Value = 350
$ws.Range("K20").Value = 2800
$ws.Range("L20").Value = 3000
$ws.Range("M20").Value = 2886
$ws.Range("P20").Value = 962
$ws.Range("D21").Value = 44460
$ws.Range("I21").Value = "Primera"
$ws.Range("J21").Value = 300
$ws.Range("K21").Value = 950
$ws.Range("L21").Value = 1000
$ws.Range("M21").Value = 975
$ws.Range("P21").Value = 325
$ws.Range("D22").Value = 44467
$ws.Range("I22").Value = "Primera"
$ws.Range("J22").Value = 250
$ws.Range("K22").Value = 800
$ws.Range("L22").Value = 900
$ws.Range("M22").Value = 850
$ws.Range("P22").Value = 283
$ws.Range("D23").Value = 44565
$ws.Range("I23").Value = "Primera"
$ws.Range("J23").Value = 250
$ws.Range("K23").Value = 3000
$ws.Range("L23").Value = 3500
$ws.Range("M23").Value = 3250
$ws.Range("P23").Value = 1083
$ws.Range("D24").Value = 44412
$ws.Range("I24").Value = "Primera"
$ws.Range("J24").Value = 300
$ws.Range("K24").Value = 2800
$ws.Range("L24").Value = 3000
$ws.Range("M24").Value = 2900
$ws.Range("P24").Value = 967
$ws.Range("D25").Value = 44483
$ws.Range("I25").Value = "Primera"
$ws.Range("J25").Value = 300
$ws.Range("K25").Value = 1000
$ws.Range("L25").Value = 1200
$ws.Range("M25").Value = 1100
$ws.Range("P25").Value = 367
$ws.Range("D26").Value = 44168
$ws.Range("I26").Value = "Primera"
$ws.Range("J26").Value = 300
$ws.Range("K26").Value = 1800
$ws.Range("L26").Value = 2000
$ws.Range("M26").Value = 1900
$ws.Range("P26").Value = 633
$ws.Range("D27").Value = 44308
$ws.Range("I27").Value = "Primera"
$ws.Range("J27").Value = 270
$ws.Range("K27").Value = 1400
$ws.Range("L27").Value = 1500
$ws.Range("M27").Value = 1450
$ws.Range("P27").Value = 483
$ws.Range("D28").Value = 44579
$ws.Range("I28").Value = "Primera"
$ws.Range("J28").Value = 300
$ws.Range("K28").Value = 3000
$ws.Range("L28").Value = 3500
$ws.Range("M28").Value = 3250
$ws.Range("P28").Value = 1083
$ws.Range("D29").Value = 44435
$ws.Range("I29").Value = "Primera"
$ws.Range("J29").Value = 270
$ws.Range("K29").Value = 1800
$ws.Range("L29").Value = 2000
$ws.Range("M29").Value = 1900
$ws.Range("P29").Value = 633
$ws.Range("D30").Value = 44350
$ws.Range("I30").Value = "Primera"
$ws.Range("J30").Value = 300
$ws.Range("K30").Value = 1800
$ws.Range("L30").Value = 2000
$ws.Range("M30").Value = 1900
$ws.Range("P30").Value = 633
$ws.Range("D31").Value = 44586
$ws.Range("I31").Value = "Primera"
$ws.Range("J31").Value = 250
$ws.Range("K31").Value = 2500
$ws.Range("L31").Value = 3000
$ws.Range("M31").Value = 2750
$ws.Range("P31").Value = 917
$ws.Range("D32").Value = 44558
$ws.Range("I32").Value = "Primera"
$ws.Range("J32").Value = 250
$ws.Range("K32").Value = 3500
$ws.Range("L32").Value = 4000
$ws.Range("M32").Value = 3750
$ws.Range("P32").Value = 1250
$ws.Range("D33").Value = 44356
$ws.Range("I33").Value = "Primera"
$ws.Range("J33").Value = 200
$ws.Range("K33").Value = 2400
$ws.Range("L33").Value = 2500
$ws.Range("M33").Value = 2450
$ws.Range("P33").Value = 817
$ws.Range("D34").Value = 44356
$ws.Range("I34").Value = "Segunda"
$ws.Range("J34").Value = 200
$ws.Range("K34").Value = 1800
$ws.Range("L34").Value = 2000
$ws.Range("M34").Value = 1900
$ws.Range("P34").Value = 633
$ws.Range("D35").Value = 44530
$ws.Range("I35").Value = "Primera"
$ws.Range("J35").Value = 300
$ws.Range("K35").Value = 1900
$ws.Range("L35").Value = 2000
$ws.Range("M35").Value = 1950
$ws.Range("P35").Value = 650
$ws.Range("D36").Value = 44487
$ws.Range("I36").Value = "Primera"
$ws.Range("J36").Value = 300
$ws.Range("K36").Value = 950
$ws.Range("L36").Value = 1000
$ws.Range("M36").Value = 975
$ws.Range("P36").Value = 325
$ws.Range("D37").Value = 44327
$ws.Range("I37").Value = "Primera"
$ws.Range("J37").Value = 200
$ws.Range("K37").Value = 1400
$ws.Range("L37").Value = 1500
$ws.Range("M37").Value = 1450
$ws.Range("P37").Value = 483
$ws.Range("D38").Value = 44327
$ws.Range("I38").Value = "Segunda"
$ws.Range("J38").Value = 250
$ws.Range("K38").Value = 1000
$ws.Range("L38").Value = 1200
$ws.Range("M38").Value = 1100
$ws.Range("P38").Value = 367
$ws.Range("D39").Value = 44364
$ws.Range("I39").Value = "Primera"
$ws.Range("J39").Value = 270
$ws.Range("K39").Value = 3400
$ws.Range("L39").Value = 3500
$ws.Range("M39").Value = 3450
$ws.Range("P39").Value = 1150
$ws.Range("D40").Value = 44376
$ws.Range("I40").Value = "Primera"
$ws.Range("J40").Value = 280
$ws.Range("K40").Value = 2400
$ws.Range("L40").Value = 2500
$ws.Range("M40").Value = 2436
$ws.Range("P40").Value = 812
$ws.Range("D41").Value = 44313
$ws.Range("I41").Value = "Primera"
$ws.Range("J41").Value = 300
$ws.Range("K41").Value = 1300
$ws.Range("L41").Value = 1500
$ws.Range("M41").Value = 1400
$ws.Range("P41").Value = 467
$ws.Range("D42").Value = 44313
$ws.Range("I42").Value = "Segunda"
$ws.Range("J42").Value = 250
$ws.Range("K42").Value = 900
$ws.Range("L42").Value = 1000
$ws.Range("M42").Value = 950
$ws.Range("P42").Value = 317